# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.419.85"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "2.626.23"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("E9").Value = "  -4.15%  "
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("D13").Value = "3.092.50"
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "62.336.74"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "2.630.57"
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.498"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("D28").Value = "0.0₃0833"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("E31").Value = "  -4.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "337.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.611"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0961"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0237"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.26%  "
